$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the last row (ZKFair / ZKF / 9636)
# ------------------------------------------------------------------
$ws.Rows.Item(12).Delete()

# ------------------------------------------------------------------
# 2. Re-label header row and introduce the two new "price" columns.
#    C (Amount) stays put; D becomes a brand new empty "Price per
#    token Now" column; the previously-empty "Price per token"
#    column slides into E and is relabelled "Price per token from
#    last check"; the old "Changes per Day" column slides into F and
#    becomes "Changes"; the old "Price" formula column slides into G.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Full name"
$ws.Range("B1").Value = "Short Name"
$ws.Range("C1").Value = "Amount"
$ws.Range("D1").Value = "Price per token Now"
$ws.Range("E1").Value = "Price per token from last check"
$ws.Range("F1").Value = "Changes"
$ws.Range("G1").Value = "Price"
$ws.Range("G1").Font.Bold = $true

# make sure the (now relocated) "last check" column starts out blank
$ws.Range("D2:D11").ClearContents()
$ws.Range("E2:E11").ClearContents()

# ------------------------------------------------------------------
# 3. Formulas
#    F = Changes = ((PriceNow last check / price now) * 1) - 1
#    G = Price   = Amount * Price per token Now
# ------------------------------------------------------------------
$ws.Range("F2").Formula = "=((E2/D2)*1)-1"
$ws.Range("F3:F11").Formula = "=((E3/D3)*1)-1"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=C$r*D$r"
}

# ------------------------------------------------------------------
# 4. Formatting
# ------------------------------------------------------------------
# Header row is taller & wraps text
$ws.Rows.Item(1).RowHeight = 30

$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").WrapText = $true

$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").WrapText = $true

# "Changes" column is shown as a percentage
$ws.Range("F2:F11").Style = "Percent"
$ws.Range("F2:F11").NumberFormat = "0.00%"

# ------------------------------------------------------------------
# 5. Selection / view state
# ------------------------------------------------------------------
$ws.Range("E13").Select()

Write-Output "done"
